$d = $word.ActiveDocument

# Append a new glossary entry paragraph after the last paragraph
# ("Especialista: ...") containing "Observador: <definition>".
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$np = $d.Paragraphs.Item($d.Paragraphs.Count)
$ir = $np.Range
$ir.Collapse(1)

# Bold label run.
$ir.InsertAfter("Observador: ")
$ir.Font.Bold = 1
$ir.Font.BoldBi = 1
$ir.Collapse(0)

# Regular definition run.
$ir.InsertAfter("Alguém que entra no sistema e não está logado.")
$ir.Font.Bold = 0
$ir.Font.BoldBi = 1
$ir.Collapse(0)

# Work around collapsed-bookmark placement landing on paragraph/body
# start when the target position sits immediately before a paragraph
# mark: type a temporary two-char placeholder so the bookmark position
# is no longer "end of paragraph minus one", add the bookmark there,
# then remove the placeholder again.
$ir.InsertAfter("XY")
$ir.Collapse(0)
$bmPos = $ir.Start - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$delRange = $d.Range($ir.Start - 2, $ir.Start)
$delRange.Delete()
